$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6924.75
$ws.Range("J62").Value = 6850
$ws.Range("L62").Value = 6850
$ws.Range("N62").Value = -8098
$ws.Range("H64").Value = 12699.667
$ws.Range("J64").Value = 14798.429
$ws.Range("L64").Value = 14798.429
$ws.Range("N64").Value = -15294.429
$ws.Range("H65").Value = 6924.75
$ws.Range("J65").Value = 6850
$ws.Range("L65").Value = 34250
$ws.Range("N65").Value = -40490
$ws.Range("H67").Value = 12699.667
$ws.Range("J67").Value = 14798.429
$ws.Range("L67").Value = 14798.429
$ws.Range("N67").Value = -16514.429
$ws.Range("H80").Value = 2253.5757
$ws.Range("I80").Value = 617.8333
$ws.Range("J80").Value = 3188.2856
$ws.Range("K80").Value = 1853.4999
$ws.Range("L80").Value = 9564.856800000001
$ws.Range("M80").Value = -855.4999
$ws.Range("N80").Value = -11560.8568
$ws.Range("H83").Value = 2253.5757
$ws.Range("I83").Value = 617.8333
$ws.Range("J83").Value = 3188.2856
$ws.Range("K83").Value = 5560.4997
$ws.Range("L83").Value = 28694.5704
$ws.Range("M83").Value = -568.4997000000003
$ws.Range("N83").Value = -38678.5704
$ws.Range("H92").Value = 364.57144
$ws.Range("I92").Value = 241.4
$ws.Range("J92").Value = 672.5
$ws.Range("K92").Value = 241.4
$ws.Range("L92").Value = 672.5
$ws.Range("M92").Value = 1006.6
$ws.Range("N92").Value = -3168.5
$ws.Range("H97").Value = 946
$ws.Range("J97").Value = 946
$ws.Range("L97").Value = 2838
$ws.Range("N97").Value = -3830
$ws.Range("H100").Value = 7281.6
$ws.Range("J100").Value = 7757.8335
$ws.Range("L100").Value = 7757.8335
$ws.Range("N100").Value = -8839.833500000001
$ws.Range("H113").Value = 12830.0625
$ws.Range("I113").Value = 7311.3335
$ws.Range("K113").Value = 7311.3335
$ws.Range("M113").Value = -4057.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2740.4375
$ws.Range("I45").Value = 2756.4666
$ws.Range("K45").Value = 2756.4666
$ws.Range("M45").Value = -2379.4666
$ws.Range("H74").Value = 8429.625
$ws.Range("I74").Value = 8739.5
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 8739.5
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -7865.5
$ws.Range("N74").Value = -9248
$ws.Range("H77").Value = 8429.625
$ws.Range("I77").Value = 8739.5
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 43697.5
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -39329.5
$ws.Range("N77").Value = -46236
$ws.Range("H123").Value = 39272.727
$ws.Range("J123").Value = 39272.727
$ws.Range("L123").Value = 39272.727
$ws.Range("N123").Value = -49072.727
$ws.Range("H132").Value = 2489.0588
$ws.Range("I132").Value = 2277.5
$ws.Range("J132").Value = 2996.8
$ws.Range("K132").Value = 6832.5
$ws.Range("L132").Value = 8990.400000000001
$ws.Range("M132").Value = -4302.5
$ws.Range("N132").Value = -14050.4
$ws.Range("H135").Value = 39769.23
$ws.Range("J135").Value = 39769.23
$ws.Range("L135").Value = 39769.23
$ws.Range("N135").Value = -49909.23
$ws.Range("H139").Value = 55500
$ws.Range("J139").Value = 55263.156
$ws.Range("L139").Value = 55263.156
$ws.Range("N139").Value = -65543.156

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19276516
$ws.Range("I86").Value = 31320096
$ws.Range("K86").Value = 31320096
$ws.Range("M86").Value = -31318973
$ws.Range("H89").Value = 19276516
$ws.Range("I89").Value = 31320096
$ws.Range("K89").Value = 156600480
$ws.Range("M89").Value = -156594864
$ws.Range("H134").Value = 2226
$ws.Range("I134").Value = 1829.8572
$ws.Range("K134").Value = 5489.571599999999
$ws.Range("M134").Value = -2954.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1365.9166
$ws.Range("I16").Value = 1165.6666
$ws.Range("K16").Value = 1165.6666
$ws.Range("M16").Value = -878.6666
$ws.Range("H113").Value = 1365.9166
$ws.Range("I113").Value = 1165.6666
$ws.Range("K113").Value = 1165.6666
$ws.Range("M113").Value = 1004.3334
$ws.Range("H132").Value = 2934.6843
$ws.Range("I132").Value = 2985.8235
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 8957.470499999999
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -6427.470499999999
$ws.Range("N132").Value = -12560
$ws.Range("H134").Value = 3099.4546
$ws.Range("I134").Value = 1718
$ws.Range("K134").Value = 5154
$ws.Range("M134").Value = -2619

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 424.25
$ws.Range("I14").Value = 424.25
$ws.Range("K14").Value = 1272.75
$ws.Range("M14").Value = -1099.75
$ws.Range("H98").Value = 1258.3043
$ws.Range("I98").Value = 1479.5834
$ws.Range("J98").Value = 1016.9091
$ws.Range("K98").Value = 4438.7502
$ws.Range("L98").Value = 3050.7273
$ws.Range("M98").Value = -2940.7502
$ws.Range("N98").Value = -6046.7273
$ws.Range("H114").Value = 861.1429000000001
$ws.Range("I114").Value = 785.6
$ws.Range("K114").Value = 2356.8
$ws.Range("M114").Value = 897.1999999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5965.6665
$ws.Range("I122").Value = 5398.636
$ws.Range("K122").Value = 16195.908
$ws.Range("M122").Value = -13745.908
$ws.Range("H126").Value = 1927.5
$ws.Range("I126").Value = 600
$ws.Range("K126").Value = 1800
$ws.Range("M126").Value = 670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2550.0417
$ws.Range("I132").Value = 2388.682
$ws.Range("K132").Value = 7166.045999999999
$ws.Range("M132").Value = -4636.045999999999
$ws.Range("H133").Value = 79950
$ws.Range("J133").Value = 79950
$ws.Range("L133").Value = 79950
$ws.Range("N133").Value = -85010
$ws.Range("H136").Value = 2989.3333
$ws.Range("I136").Value = 3225.25
$ws.Range("K136").Value = 9675.75
$ws.Range("M136").Value = -7125.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 33856.145
$ws.Range("I51").Value = 18500
$ws.Range("K51").Value = 18500
$ws.Range("M51").Value = -17990
$ws.Range("H122").Value = 2508.7856
$ws.Range("I122").Value = 1966
$ws.Range("J122").Value = 4499
$ws.Range("K122").Value = 5898
$ws.Range("L122").Value = 13497
$ws.Range("M122").Value = -3448
$ws.Range("N122").Value = -18397
$ws.Range("H132").Value = 4750.7
$ws.Range("I132").Value = 2446.3333
$ws.Range("K132").Value = 7338.999899999999
$ws.Range("M132").Value = -4808.999899999999
